$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Right single quotation mark (U+2019), used by the original author's text.
$rsquo = [char]0x2019

# Replace the CONCATENATE formulas in D1:D9 with literal, zero-padded
# numbered text values (typed over the old formula results).
$ws.Range("D1").Value = "01 Copy of the legal title report "
$ws.Range("D2").Value = "02 Details of encumbrances "
$ws.Range("D3").Value = "03 Copy of Layout Approval (in case of layout) "
$ws.Range("D4").Value = "04 Building Plan Approval / NA Order for plotted development "
$ws.Range("D5").Value = "05 Commencement Certificates / NA Order for plotted development "
$ws.Range("D6").Value = "06 Declaration about Commencement Certificate "
$ws.Range("D7").Value = "07 Declaration in FORM B "
$ws.Range("D8").Value = "08 Architect" + $rsquo + "s Certificate of Percentage of Completion of Work (Form 1)"
$ws.Range("D9").Value = "09 Engineer" + $rsquo + "s Certificate on Cost Incurred on Project (Form 2)"

# Re-establish the CONCATENATE shared formula across the surviving rows
# (D10:D19) now that D1:D9 hold literal values instead of formulas.
$ws.Range("D10:D19").Formula = '=CONCATENATE(A10," ",B10)'

# Give column D the same width treatment as column B.
$ws.Columns("D").ColumnWidth = $ws.Columns("B").ColumnWidth

# Move the active selection to J11 (was J19).
$ws.Range("J11").Select()
